$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# Update header row
$ws.Range("B1").Value = "Bid"
$ws.Range("C1").Value = "Ask"
$ws.Range("D1").Value = "Last"

# Update data row
$ws.Range("B2").Value = "-"

# Remove the now-unused column E (High/Low -> removed) to shrink dimension to A1:D2
$ws.Range("E1:E2").Delete()
